$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1 / D1 get a top+bottom border (matches border id 4 used elsewhere)
$ws1.Range("C1").Borders.Item(3).LineStyle = 1   # xlEdgeTop
$ws1.Range("C1").Borders.Item(4).LineStyle = 1   # xlEdgeBottom

# D1 gets a top+bottom+right border (matches border id 5)
$ws1.Range("D1").Borders.Item(3).LineStyle = 1   # xlEdgeTop
$ws1.Range("D1").Borders.Item(4).LineStyle = 1   # xlEdgeBottom
$ws1.Range("D1").Borders.Item(2).LineStyle = 1   # xlEdgeRight

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2.Range("C1").Borders.Item(3).LineStyle = 1
$ws2.Range("C1").Borders.Item(4).LineStyle = 1

$ws2.Range("D1").Borders.Item(3).LineStyle = 1
$ws2.Range("D1").Borders.Item(4).LineStyle = 1
$ws2.Range("D1").Borders.Item(2).LineStyle = 1

$ws2.Range("F1").Borders.Item(3).LineStyle = 1
$ws2.Range("F1").Borders.Item(4).LineStyle = 1

$ws2.Range("G1").Borders.Item(3).LineStyle = 1
$ws2.Range("G1").Borders.Item(4).LineStyle = 1
$ws2.Range("G1").Borders.Item(2).LineStyle = 1

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 becomes an empty cell (was an empty inline string cell, now cleared entirely)
$ws2.Range("G5").ClearContents()
